# Edit script for other_transport.xlsx
# Fixes misaligned annual data rows (1958-1974): consolidates values that had
# drifted across rows 2-33 back into their correct year row (rows 2-18), adds
# several previously-missing nonmilitary-industry data points, and cleans up
# the now-empty trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 7.5
$ws.Range("L3").Value = 2.8
$ws.Range("F4").Value = 2.06
$ws.Range("K4").Value = 9.9
$ws.Range("L4").Value = 2.1
$ws.Range("M4").Value = 7
$ws.Range("B5").Value = 30.3
$ws.Range("D5").Value = 24.3
$ws.Range("F5").Value = 2.13
$ws.Range("G5").Value = 39.299999999999997
$ws.Range("H5").Value = 1.8
$ws.Range("I5").Value = 7.8
$ws.Range("J5").Value = 3.9
$ws.Range("K5").Value = 6.9
$ws.Range("L5").Value = 2.2000000000000002
$ws.Range("M5").Value = 3.8
$ws.Range("B6").Value = 35.4
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 28.6
$ws.Range("F6").Value = 2.15
$ws.Range("G6").Value = 40.1
$ws.Range("H6").Value = 2.5
$ws.Range("I6").Value = 8.1999999999999993
$ws.Range("J6").Value = 5.3
$ws.Range("K6").Value = 7.4
$ws.Range("L6").Value = 2.7
$ws.Range("M6").Value = 3.5
$ws.Range("B7").Value = 42.6
$ws.Range("C7").Value = 4.9000000000000004
$ws.Range("D7").Value = 34.799999999999997
$ws.Range("E7").Value = 91.84
$ws.Range("F7").Value = 2.2400000000000002
$ws.Range("G7").Value = 41
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 5.0999999999999996
$ws.Range("K7").Value = 7.1
$ws.Range("L7").Value = 2.9
$ws.Range("M7").Value = 2.8
$ws.Range("B8").Value = 50.7
$ws.Range("C8").Value = 5.7
$ws.Range("D8").Value = 41.8
$ws.Range("E8").Value = 93.89
$ws.Range("F8").Value = 2.29
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 8.3000000000000007
$ws.Range("J8").Value = 6.4
$ws.Range("K8").Value = 7.5
$ws.Range("L8").Value = 3.6
$ws.Range("M8").Value = 2.6
$ws.Range("B9").Value = 57.3
$ws.Range("C9").Value = 6.1
$ws.Range("D9").Value = 47.1
$ws.Range("E9").Value = 93.09
$ws.Range("F9").Value = 2.31
$ws.Range("G9").Value = 40.299999999999997
$ws.Range("H9").Value = 2.9
$ws.Range("I9").Value = 8.1999999999999993
$ws.Range("J9").Value = 6.5
$ws.Range("K9").Value = 8
$ws.Range("L9").Value = 4.2
$ws.Range("M9").Value = 2.4
$ws.Range("B10").Value = 63.4
$ws.Range("C10").Value = 7.3
$ws.Range("D10").Value = 52.2
$ws.Range("E10").Value = 95.52
$ws.Range("F10").Value = 2.4
$ws.Range("G10").Value = 39.799999999999997
$ws.Range("H10").Value = 2.7
$ws.Range("I10").Value = 8.4
$ws.Range("J10").Value = 7
$ws.Range("K10").Value = 10.4
$ws.Range("L10").Value = 5.4
$ws.Range("M10").Value = 3
$ws.Range("B11").Value = 67.2
$ws.Range("C11").Value = 8.5
$ws.Range("D11").Value = 55.3
$ws.Range("E11").Value = 102.11
$ws.Range("F11").Value = 2.54
$ws.Range("G11").Value = 40.200000000000003
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = 9.6
$ws.Range("J11").Value = 7.7
$ws.Range("K11").Value = 8.6999999999999993
$ws.Range("L11").Value = 4.9000000000000004
$ws.Range("M11").Value = 1.9
$ws.Range("B12").Value = 84.6
$ws.Range("C12").Value = 11.1
$ws.Range("D12").Value = 71
$ws.Range("E12").Value = 111.08
$ws.Range("F12").Value = 2.77
$ws.Range("G12").Value = 40.1
$ws.Range("I12").Value = 9.4
$ws.Range("J12").Value = 8.4
$ws.Range("K12").Value = 9
$ws.Range("L12").Value = 5.4
$ws.Range("M12").Value = 1.5
$ws.Range("B13").Value = 104.7
$ws.Range("C13").Value = 13.7
$ws.Range("D13").Value = 87
$ws.Range("E13").Value = 115.53
$ws.Range("F13").Value = 2.91
$ws.Range("G13").Value = 39.700000000000003
$ws.Range("H13").Value = 2.8
$ws.Range("I13").Value = 8.6
$ws.Range("K13").Value = 10
$ws.Range("L13").Value = 5.6
$ws.Range("M13").Value = 2.2999999999999998
$ws.Range("B14").Value = 111.4
$ws.Range("C14").Value = 15.2
$ws.Range("D14").Value = 91.5
$ws.Range("E14").Value = 118.81
$ws.Range("F14").Value = 3.07
$ws.Range("G14").Value = 38.700000000000003
$ws.Range("H14").Value = 2.2000000000000002
$ws.Range("I14").Value = 7.9
$ws.Range("J14").Value = 6.3
$ws.Range("K14").Value = 8.5
$ws.Range("L14").Value = 4.2
$ws.Range("M14").Value = 2.5
$ws.Range("B15").Value = 131.9
$ws.Range("C15").Value = 17.3
$ws.Range("D15").Value = 108.6
$ws.Range("E15").Value = 129.75
$ws.Range("F15").Value = 3.31
$ws.Range("G15").Value = 39.200000000000003
$ws.Range("H15").Value = 2.6
$ws.Range("I15").Value = 8.1999999999999993
$ws.Range("J15").Value = 6.4
$ws.Range("K15").Value = 6.5
$ws.Range("L15").Value = 3.6
$ws.Range("M15").Value = 1.5
$ws.Range("B16").Value = 162.9
$ws.Range("C16").Value = 22.8
$ws.Range("D16").Value = 135.5
$ws.Range("E16").Value = 138.94999999999999
$ws.Range("F16").Value = 3.5
$ws.Range("G16").Value = 39.700000000000003
$ws.Range("H16").Value = 3.1
$ws.Range("I16").Value = 8.8000000000000007
$ws.Range("J16").Value = 7.4
$ws.Range("K16").Value = 7.8
$ws.Range("L16").Value = 4.9000000000000004
$ws.Range("M16").Value = 1.4
$ws.Range("B17").Value = 167.7
$ws.Range("C17").Value = 26.6
$ws.Range("D17").Value = 137.30000000000001
$ws.Range("E17").Value = 142.80000000000001
$ws.Range("F17").Value = 3.69
$ws.Range("G17").Value = 38.700000000000003
$ws.Range("I17").Value = 8
$ws.Range("J17").Value = 7.1
$ws.Range("K17").Value = 9.6
$ws.Range("M17").Value = 2.4
$ws.Range("B18").Value = 130.6
$ws.Range("C18").Value = 22.6
$ws.Range("D18").Value = 103.4
$ws.Range("E18").Value = 148.6
$ws.Range("F18").Value = 3.88
$ws.Range("G18").Value = 38.299999999999997
$ws.Range("H18").Value = 2.2000000000000002
$ws.Range("I18").Value = 7
$ws.Range("J18").Value = 5.6
$ws.Range("K18").Value = 10.6
$ws.Range("M18").Value = 4.9000000000000004

# Row 10 (1966) column M picks up the 2-decimal number format that used to
# live on the old M18 cell; M18 itself reverts to the default/general style.
$ws.Range("M18").Style = "Normal"
$ws.Range("M10").NumberFormat = "0.00"

# The trailing rows (19-33) held the stray/duplicate values that have now
# been folded back into rows 2-18 above, so clear them out entirely.
$ws.Range("A19:M33").ClearContents()

# Restore the selection that was active when the workbook was last saved.
$ws.Range("B10").Select()
